$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the additional HMC variations, right after the
# existing HMCv12 row (row 15), pushing the WW/TWW/ELW/AMW rows down by 2.
$ws.Rows("16:17").Insert()

# Match the row height used by the rest of the data rows in the table.
$ws.Rows("16:17").RowHeight = 18.75

# Row 16 - HMCv13
$ws.Cells.Item(16, 2).Value = "HMCv13"
$ws.Cells.Item(16, 3).Value = "HMC"
$ws.Cells.Item(16, 4).Value = 677
$ws.Cells.Item(16, 5).Value = "none"
$ws.Cells.Item(16, 6).Value = "water"
$ws.Cells.Item(16, 7).Value = "water"
$ws.Cells.Item(16, 10).Value = 67.2
$ws.Cells.Item(16, 11).Value = 32.8

# Row 17 - hmcV14
$ws.Cells.Item(17, 2).Value = "hmcV14"
$ws.Cells.Item(17, 3).Value = "HMC"
$ws.Cells.Item(17, 4).Value = 1230
$ws.Cells.Item(17, 5).Value = "iron magma"
$ws.Cells.Item(17, 6).Value = "carbon dioxide"
$ws.Cells.Item(17, 7).Value = "carbon dioxide"
$ws.Cells.Item(17, 8).Value = "nitrogen"
$ws.Cells.Item(17, 9).Value = "sulphur dioxide"
$ws.Cells.Item(17, 10).Value = 65.9
$ws.Cells.Item(17, 11).Value = 34.1

# Widen column E (it was previously auto-best-fit) to roughly match the
# author's manual resize.
$ws.Columns("E").ColumnWidth = 17.75

# Leave the selection where the author left it after the edit.
$ws.Range("L17").Select()
